# Regenerate the experiment task-order sheets: new random filenames/timestamps
# for each task-order sheet, and a reshuffled sheet tab order.
# (commit: "Created experiment order generation script")

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) GNG_TO sheet: rename + refresh its go/GNG task-order values (5 rows incl. header)
# ---------------------------------------------------------------------------
$sheetGNG = $wb.Worksheets.Item("GNG_TO-16512555138297467")
$sheetGNG.Name = "GNG_TO-16515889433005884"
$sheetGNG.Range("B2").Value = "go_stims-16515889432625701.csv"
$sheetGNG.Range("B3").Value = "GNG_stims-1651588943286016.csv"
$sheetGNG.Range("B4").Value = "go_stims-1651588943287041.csv"
$sheetGNG.Range("B5").Value = "GNG_stims-16515889432995763.csv"

# ---------------------------------------------------------------------------
# 2) NB_TO sheet: rename + refresh its OB/TB/ZB-match task-order values (10 rows incl. header)
# ---------------------------------------------------------------------------
$sheetNB = $wb.Worksheets.Item("NB_TO-16512555155708125")
$sheetNB.Name = "NB_TO-165158894312961"
$sheetNB.Range("B2").Value = "ZB-match_1-1651588940298163.csv"
$sheetNB.Range("B3").Value = "ZB-match_9-16515889401308095.csv"
$sheetNB.Range("B4").Value = "OB-16515889410189502.csv"
$sheetNB.Range("B5").Value = "TB-16515889419243107.csv"
$sheetNB.Range("B6").Value = "TB-1651588943102839.csv"
$sheetNB.Range("B7").Value = "OB-16515889403463159.csv"
$sheetNB.Range("B8").Value = "TB-16515889415524783.csv"
$sheetNB.Range("B9").Value = "ZB-match_7-16515889400746477.csv"
$sheetNB.Range("B10").Value = "OB-16515889407075605.csv"

# ---------------------------------------------------------------------------
# 3) RS_TO sheet: rename + swap the resting-state order (eyes open / eyes closed)
# ---------------------------------------------------------------------------
$sheetRS = $wb.Worksheets.Item("RS_TO-16512555155718138")
$sheetRS.Name = "RS_TO-16515889433015823"
$sheetRS.Range("B2").Value = "eyes closed"
$sheetRS.Range("B3").Value = "eyes open"

# ---------------------------------------------------------------------------
# 4) TOL_TO sheet: rename + refresh its MM/ZM task-order values (7 rows incl. header)
# ---------------------------------------------------------------------------
$sheetTOL = $wb.Worksheets.Item("TOL_TO-16512555156348152")
$sheetTOL.Name = "TOL_TO-16515889431924946"
$sheetTOL.Range("B2").Value = "MM_stims-1651588943161521.csv"
$sheetTOL.Range("B3").Value = "ZM_stims-1651588943137244.csv"
$sheetTOL.Range("B4").Value = "MM_stims-16515889431771543.csv"
$sheetTOL.Range("B5").Value = "ZM_stims-1651588943162495.csv"
$sheetTOL.Range("B6").Value = "MM_stims-16515889431924946.csv"
$sheetTOL.Range("B7").Value = "ZM_stims-16515889431781235.csv"

# ---------------------------------------------------------------------------
# 5) vSAT_TO sheet: rename + refresh its vSAT/SAT task-order values (5 rows incl. header)
# ---------------------------------------------------------------------------
$sheetVSAT = $wb.Worksheets.Item("vSAT_TO-16512555157128134")
$sheetVSAT.Name = "vSAT_TO-16515889432568698"
$sheetVSAT.Range("B2").Value = "vSAT_stims-16515889432259598.csv"
$sheetVSAT.Range("B3").Value = "SAT_stims-16515889431971092.csv"
$sheetVSAT.Range("B4").Value = "SAT_stims-16515889432090235.csv"
$sheetVSAT.Range("B5").Value = "vSAT_stims-16515889432407753.csv"

# ---------------------------------------------------------------------------
# 6) Reorder the sheet tabs to: NB_TO, TOL_TO, vSAT_TO, GNG_TO, RS_TO
#    (re-fetch each sheet by its new stable name right before moving it, so
#    the move always acts on current worksheet positions)
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("NB_TO-165158894312961").Move($wb.Worksheets.Item("GNG_TO-16515889433005884"))
$wb.Worksheets.Item("TOL_TO-16515889431924946").Move($null, $wb.Worksheets.Item("NB_TO-165158894312961"))
$wb.Worksheets.Item("vSAT_TO-16515889432568698").Move($null, $wb.Worksheets.Item("TOL_TO-16515889431924946"))
$wb.Worksheets.Item("GNG_TO-16515889433005884").Move($null, $wb.Worksheets.Item("vSAT_TO-16515889432568698"))
$wb.Worksheets.Item("RS_TO-16515889433015823").Move($null, $wb.Worksheets.Item("GNG_TO-16515889433005884"))

foreach ($ws in $wb.Worksheets) {
    Write-Output $ws.Name
}
